$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha (D) and Volumen (M) values between row 2 and row 4
$ws.Range("D2").Value = 44253
$ws.Range("M2").Value = 12

$ws.Range("D4").Value = 44993
$ws.Range("M4").Value = 14
